# Edit script: applies the "compareExtended" / SOAP / image-diff / desktop typeKeys / web clickAll
# showcase command additions to the #system reference sheet, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$vals_G = New-Object 'object[,]' 94,1
$vals_G[0,0] = "assertAttribute(locator,attribute,expected)"
$vals_G[1,0] = "assertChecked(name)"
$vals_G[2,0] = "assertDisabled(name)"
$vals_G[3,0] = "assertElementPresent(name)"
$vals_G[4,0] = "assertEnabled(name)"
$vals_G[5,0] = "assertHierCells(matchBy,column,expected,nestedOnly)"
$vals_G[6,0] = "assertHierRow(matchBy,expected)"
$vals_G[7,0] = "assertListCount(count)"
$vals_G[8,0] = "assertLocatorNotPresent(locator)"
$vals_G[9,0] = "assertLocatorPresent(locator)"
$vals_G[10,0] = "assertMenuEnabled(menu)"
$vals_G[11,0] = "assertModalDialogNotPresent()"
$vals_G[12,0] = "assertModalDialogPresent()"
$vals_G[13,0] = "assertModalDialogTitle(title)"
$vals_G[14,0] = "assertModalDialogTitleByLocator(locator,title)"
$vals_G[15,0] = "assertNotChecked(name)"
$vals_G[16,0] = "assertSelected(name,text)"
$vals_G[17,0] = "assertTableCell(row,column,contains)"
$vals_G[18,0] = "assertTableColumnContains(column,contains)"
$vals_G[19,0] = "assertTableContains(contains)"
$vals_G[20,0] = "assertTableRowContains(row,contains)"
$vals_G[21,0] = "assertText(name,expected)"
$vals_G[22,0] = "assertWindowTitleContains(contains)"
$vals_G[23,0] = "clear(locator)"
$vals_G[24,0] = "clearCombo(name)"
$vals_G[25,0] = "clearModalDialog(var,button)"
$vals_G[26,0] = "clearTextArea(name)"
$vals_G[27,0] = "clearTextBox(name)"
$vals_G[28,0] = "clickButton(name)"
$vals_G[29,0] = "clickByLocator(locator)"
$vals_G[30,0] = "clickCheckBox(name)"
$vals_G[31,0] = "clickElementOffset(name,xOffset,yOffset)"
$vals_G[32,0] = "clickExplorerBar(group,item)"
$vals_G[33,0] = "clickFirstMatchRow(nameValues)"
$vals_G[34,0] = "clickFirstMatchedList(contains)"
$vals_G[35,0] = "clickIcon(label)"
$vals_G[36,0] = "clickList(row)"
$vals_G[37,0] = "clickMenu(menu)"
$vals_G[38,0] = "clickOffset(locator,xOffset,yOffset)"
$vals_G[39,0] = "clickRadio(name)"
$vals_G[40,0] = "clickTab(group,name)"
$vals_G[41,0] = "clickTableCell(row,column)"
$vals_G[42,0] = "clickTableRow(row)"
$vals_G[43,0] = "clickTextPane(name,criteria)"
$vals_G[44,0] = "clickTextPaneRow(var,index)"
$vals_G[45,0] = "closeApplication()"
$vals_G[46,0] = "collapseHierTable()"
$vals_G[47,0] = "editCurrentRow(nameValues)"
$vals_G[48,0] = "editHierCells(var,matchBy,nameValues)"
$vals_G[49,0] = "editTableCells(row,nameValues)"
$vals_G[50,0] = "getRowCount(var)"
$vals_G[51,0] = "hideExplorerBar()"
$vals_G[52,0] = "login(form,username,password)"
$vals_G[53,0] = "maximize()"
$vals_G[54,0] = "minimize()"
$vals_G[55,0] = "resize(width,height)"
$vals_G[56,0] = "saveAllTableRows(var)"
$vals_G[57,0] = "saveAttributeByLocator(var,locator,attribute)"
$vals_G[58,0] = "saveElementCount(var,name)"
$vals_G[59,0] = "saveFirstListData(var,contains)"
$vals_G[60,0] = "saveFirstMatchedListIndex(var,contains)"
$vals_G[61,0] = "saveHierCells(var,matchBy,column,nestedOnly)"
$vals_G[62,0] = "saveHierRow(var,matchBy)"
$vals_G[63,0] = "saveListData(var,contains)"
$vals_G[64,0] = "saveLocatorCount(var,locator)"
$vals_G[65,0] = "saveModalDialogText(var)"
$vals_G[66,0] = "saveModalDialogTextByLocator(var,locator)"
$vals_G[67,0] = "saveProcessId(var,locator)"
$vals_G[68,0] = "saveRowCount(var)"
$vals_G[69,0] = "saveTableRows(var,contains)"
$vals_G[70,0] = "saveTableRowsRange(var,beginRow,endRow)"
$vals_G[71,0] = "saveText(var,name)"
$vals_G[72,0] = "saveTextByLocator(var,locator)"
$vals_G[73,0] = "saveTextPane(var,name,criteria)"
$vals_G[74,0] = "saveWindowTitle(var)"
$vals_G[75,0] = "scanTable(var,name)"
$vals_G[76,0] = "selectCombo(name,text)"
$vals_G[77,0] = "sendKeysToTextBox(name,text1,text2,text3,text4)"
$vals_G[78,0] = "showExplorerBar()"
$vals_G[79,0] = "toggleExplorerBar()"
$vals_G[80,0] = "typeAppendTextArea(name,text1,text2,text3,text4)"
$vals_G[81,0] = "typeAppendTextBox(name,text1,text2,text3,text4)"
$vals_G[82,0] = "typeByLocator(locator,text)"
$vals_G[83,0] = "typeKeys(os,keystrokes)"
$vals_G[84,0] = "typeTextArea(name,text1,text2,text3,text4)"
$vals_G[85,0] = "typeTextBox(name,text1,text2,text3,text4)"
$vals_G[86,0] = "useApp(appId)"
$vals_G[87,0] = "useForm(formName)"
$vals_G[88,0] = "useHierTable(var,name)"
$vals_G[89,0] = "useList(var,name)"
$vals_G[90,0] = "useTable(var,name)"
$vals_G[91,0] = "useTableRow(var,row)"
$vals_G[92,0] = "waitFor(name,maxWaitMs)"
$vals_G[93,0] = "waitForLocator(locator,maxWaitMs)"
$ws.Range("G2:G95").Value = $vals_G

$vals_J = New-Object 'object[,]' 6,1
$vals_J[0,0] = "colorbit(source,bit,saveTo)"
$vals_J[1,0] = "compare(baseline,actual)"
$vals_J[2,0] = "convert(source,format,saveTo)"
$vals_J[3,0] = "crop(image,dimension,saveTo)"
$vals_J[4,0] = "resize(image,width,height,saveTo)"
$vals_J[5,0] = "saveDiff(var,baseline,actual)"
$ws.Range("J2:J7").Value = $vals_J

$vals_M = New-Object 'object[,]' 16,1
$vals_M[0,0] = "addOrReplace(json,jsonpath,input,var)"
$vals_M[1,0] = "assertCorrectness(json,schema)"
$vals_M[2,0] = "assertElementCount(json,jsonpath,count)"
$vals_M[3,0] = "assertElementNotPresent(json,jsonpath)"
$vals_M[4,0] = "assertElementPresent(json,jsonpath)"
$vals_M[5,0] = "assertEqual(expected,actual)"
$vals_M[6,0] = "assertValue(json,jsonpath,expected)"
$vals_M[7,0] = "assertValues(json,jsonpath,array,exactOrder)"
$vals_M[8,0] = "assertWellformed(json)"
$vals_M[9,0] = "beautify(json,var)"
$vals_M[10,0] = "compact(var,json,removeEmpty)"
$vals_M[11,0] = "fromCsv(csv,header,jsonFile)"
$vals_M[12,0] = "minify(json,var)"
$vals_M[13,0] = "storeCount(json,jsonpath,var)"
$vals_M[14,0] = "storeValue(json,jsonpath,var)"
$vals_M[15,0] = "storeValues(json,jsonpath,var)"
$ws.Range("M2:M17").Value = $vals_M

$vals_Y = New-Object 'object[,]' 128,1
$vals_Y[0,0] = "assertAndClick(locator,label)"
$vals_Y[1,0] = "assertAttribute(locator,attrName,value)"
$vals_Y[2,0] = "assertAttributeContains(locator,attrName,contains)"
$vals_Y[3,0] = "assertAttributeNotContains(locator,attrName,contains)"
$vals_Y[4,0] = "assertAttributeNotPresent(locator,attrName)"
$vals_Y[5,0] = "assertAttributePresent(locator,attrName)"
$vals_Y[6,0] = "assertChecked(locator)"
$vals_Y[7,0] = "assertContainCount(locator,text,count)"
$vals_Y[8,0] = "assertCssNotPresent(locator,property)"
$vals_Y[9,0] = "assertCssPresent(locator,property,value)"
$vals_Y[10,0] = "assertElementByAttributes(nameValues)"
$vals_Y[11,0] = "assertElementByText(locator,text)"
$vals_Y[12,0] = "assertElementCount(locator,count)"
$vals_Y[13,0] = "assertElementNotPresent(locator)"
$vals_Y[14,0] = "assertElementPresent(locator)"
$vals_Y[15,0] = "assertElementsPresent(prefix)"
$vals_Y[16,0] = "assertFocus(locator)"
$vals_Y[17,0] = "assertFrameCount(count)"
$vals_Y[18,0] = "assertFramePresent(frameName)"
$vals_Y[19,0] = "assertIECompatMode()"
$vals_Y[20,0] = "assertIENativeMode()"
$vals_Y[21,0] = "assertLinkByLabel(label)"
$vals_Y[22,0] = "assertNotChecked(locator)"
$vals_Y[23,0] = "assertNotFocus(locator)"
$vals_Y[24,0] = "assertNotText(locator,text)"
$vals_Y[25,0] = "assertNotVisible(locator)"
$vals_Y[26,0] = "assertOneMatch(locator)"
$vals_Y[27,0] = "assertScrollbarHNotPresent(locator)"
$vals_Y[28,0] = "assertScrollbarHPresent(locator)"
$vals_Y[29,0] = "assertScrollbarVNotPresent(locator)"
$vals_Y[30,0] = "assertScrollbarVPresent(locator)"
$vals_Y[31,0] = "assertTable(locator,row,column,text)"
$vals_Y[32,0] = "assertText(locator,text)"
$vals_Y[33,0] = "assertTextContains(locator,text)"
$vals_Y[34,0] = "assertTextCount(locator,text,count)"
$vals_Y[35,0] = "assertTextList(locator,list,ignoreOrder)"
$vals_Y[36,0] = "assertTextMatches(text,minMatch,scrollTo)"
$vals_Y[37,0] = "assertTextNotContains(locator,text)"
$vals_Y[38,0] = "assertTextNotPresent(text)"
$vals_Y[39,0] = "assertTextOrder(locator,descending)"
$vals_Y[40,0] = "assertTextPresent(text)"
$vals_Y[41,0] = "assertTitle(text)"
$vals_Y[42,0] = "assertValue(locator,value)"
$vals_Y[43,0] = "assertValueOrder(locator,descending)"
$vals_Y[44,0] = "assertVisible(locator)"
$vals_Y[45,0] = "checkAll(locator)"
$vals_Y[46,0] = "clearLocalStorage()"
$vals_Y[47,0] = "click(locator)"
$vals_Y[48,0] = "clickAll(locator)"
$vals_Y[49,0] = "clickAndWait(locator,waitMs)"
$vals_Y[50,0] = "clickByLabel(label)"
$vals_Y[51,0] = "clickByLabelAndWait(label,waitMs)"
$vals_Y[52,0] = "clickOffset(locator,x,y)"
$vals_Y[53,0] = "clickWithKeys(locator,keys)"
$vals_Y[54,0] = "close()"
$vals_Y[55,0] = "closeAll()"
$vals_Y[56,0] = "deselect(locator,text)"
$vals_Y[57,0] = "deselectMulti(locator,array)"
$vals_Y[58,0] = "dismissInvalidCert()"
$vals_Y[59,0] = "dismissInvalidCertPopup()"
$vals_Y[60,0] = "doubleClick(locator)"
$vals_Y[61,0] = "doubleClickAndWait(locator,waitMs)"
$vals_Y[62,0] = "doubleClickByLabel(label)"
$vals_Y[63,0] = "doubleClickByLabelAndWait(label,waitMs)"
$vals_Y[64,0] = "dragAndDrop(fromLocator,toLocator)"
$vals_Y[65,0] = "dragTo(fromLocator,xOffset,yOffset)"
$vals_Y[66,0] = "editLocalStorage(key,value)"
$vals_Y[67,0] = "executeScript(var,script)"
$vals_Y[68,0] = "focus(locator)"
$vals_Y[69,0] = "goBack()"
$vals_Y[70,0] = "goBackAndWait()"
$vals_Y[71,0] = "maximizeWindow()"
$vals_Y[72,0] = "mouseOver(locator)"
$vals_Y[73,0] = "open(url)"
$vals_Y[74,0] = "openAndWait(url,waitMs)"
$vals_Y[75,0] = "openHttpBasic(url,username,password)"
$vals_Y[76,0] = "openIgnoreTimeout(url)"
$vals_Y[77,0] = "refresh()"
$vals_Y[78,0] = "refreshAndWait()"
$vals_Y[79,0] = "resizeWindow(width,height)"
$vals_Y[80,0] = "rightClick(locator)"
$vals_Y[81,0] = "saveAllWindowIds(var)"
$vals_Y[82,0] = "saveAllWindowNames(var)"
$vals_Y[83,0] = "saveAttribute(var,locator,attrName)"
$vals_Y[84,0] = "saveAttributeList(var,locator,attrName)"
$vals_Y[85,0] = "saveCount(var,locator)"
$vals_Y[86,0] = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$vals_Y[87,0] = "saveElement(var,locator)"
$vals_Y[88,0] = "saveElements(var,locator)"
$vals_Y[89,0] = "saveLocalStorage(var,key)"
$vals_Y[90,0] = "saveLocation(var)"
$vals_Y[91,0] = "savePageAs(var,sessionIdName,url)"
$vals_Y[92,0] = "savePageAsFile(sessionIdName,url,file)"
$vals_Y[93,0] = "saveTableAsCsv(locator,nextPageLocator,file)"
$vals_Y[94,0] = "saveText(var,locator)"
$vals_Y[95,0] = "saveTextArray(var,locator)"
$vals_Y[96,0] = "saveTextSubstringAfter(var,locator,delim)"
$vals_Y[97,0] = "saveTextSubstringBefore(var,locator,delim)"
$vals_Y[98,0] = "saveTextSubstringBetween(var,locator,start,end)"
$vals_Y[99,0] = "saveValue(var,locator)"
$vals_Y[100,0] = "saveValues(var,locator)"
$vals_Y[101,0] = "scrollElement(locator,xOffset,yOffset)"
$vals_Y[102,0] = "scrollLeft(locator,pixel)"
$vals_Y[103,0] = "scrollPage(xOffset,yOffset)"
$vals_Y[104,0] = "scrollRight(locator,pixel)"
$vals_Y[105,0] = "scrollTo(locator)"
$vals_Y[106,0] = "select(locator,text)"
$vals_Y[107,0] = "selectFrame(locator)"
$vals_Y[108,0] = "selectMulti(locator,array)"
$vals_Y[109,0] = "selectMultiOptions(locator)"
$vals_Y[110,0] = "selectText(locator)"
$vals_Y[111,0] = "selectWindow(winId)"
$vals_Y[112,0] = "selectWindowAndWait(winId,waitMs)"
$vals_Y[113,0] = "selectWindowByIndex(index)"
$vals_Y[114,0] = "selectWindowByIndexAndWait(index,waitMs)"
$vals_Y[115,0] = "toggleSelections(locator)"
$vals_Y[116,0] = "type(locator,value)"
$vals_Y[117,0] = "typeKeys(locator,value)"
$vals_Y[118,0] = "uncheckAll(locator)"
$vals_Y[119,0] = "unselectAllText()"
$vals_Y[120,0] = "upload(fieldLocator,file)"
$vals_Y[121,0] = "verifyContainText(locator,text)"
$vals_Y[122,0] = "verifyText(locator,text)"
$vals_Y[123,0] = "wait(waitMs)"
$vals_Y[124,0] = "waitForElementPresent(locator)"
$vals_Y[125,0] = "waitForPopUp(winId,waitMs)"
$vals_Y[126,0] = "waitForTextPresent(text)"
$vals_Y[127,0] = "waitForTitle(text)"
$ws.Range("Y2:Y129").Value = $vals_Y

$vals_AD = New-Object 'object[,]' 26,1
$vals_AD[0,0] = "append(xml,xpath,content,var)"
$vals_AD[1,0] = "assertCorrectness(xml,schema)"
$vals_AD[2,0] = "assertElementCount(xml,xpath,count)"
$vals_AD[3,0] = "assertElementNotPresent(xml,xpath)"
$vals_AD[4,0] = "assertElementPresent(xml,xpath)"
$vals_AD[5,0] = "assertSoap(wsdl,xml)"
$vals_AD[6,0] = "assertSoapFaultCode(expected,xml)"
$vals_AD[7,0] = "assertSoapFaultString(expected,xml)"
$vals_AD[8,0] = "assertValue(xml,xpath,expected)"
$vals_AD[9,0] = "assertValues(xml,xpath,array,exactOrder)"
$vals_AD[10,0] = "assertWellformed(xml)"
$vals_AD[11,0] = "beautify(xml,var)"
$vals_AD[12,0] = "clear(xml,xpath,var)"
$vals_AD[13,0] = "delete(xml,xpath,var)"
$vals_AD[14,0] = "insertAfter(xml,xpath,content,var)"
$vals_AD[15,0] = "insertBefore(xml,xpath,content,var)"
$vals_AD[16,0] = "minify(xml,var)"
$vals_AD[17,0] = "prepend(xml,xpath,content,var)"
$vals_AD[18,0] = "replace(xml,xpath,content,var)"
$vals_AD[19,0] = "replaceIn(xml,xpath,content,var)"
$vals_AD[20,0] = "storeCount(xml,xpath,var)"
$vals_AD[21,0] = "storeSoapFaultCode(var,xml)"
$vals_AD[22,0] = "storeSoapFaultDetail(var,xml)"
$vals_AD[23,0] = "storeSoapFaultString(var,xml)"
$vals_AD[24,0] = "storeValue(xml,xpath,var)"
$vals_AD[25,0] = "storeValues(xml,xpath,var)"
$ws.Range("AD2:AD27").Value = $vals_AD

# Update defined-name ranges that grew because of the newly inserted rows
$nameUpdates = @{
    "desktop" = "`$G`$2:`$G`$95"
    "image" = "`$J`$2:`$J`$7"
    "json" = "`$M`$2:`$M`$17"
    "web" = "`$Y`$2:`$Y`$129"
    "xml" = "`$AD`$2:`$AD`$27"
}
foreach ($n in $wb.Names) {
    if ($nameUpdates.ContainsKey($n.Name)) {
        $n.RefersTo = "='#system'!" + $nameUpdates[$n.Name]
    }
}
